$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new trade row (row 7)
$ws.Range("A7").Value = 42649.64471064815
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = 10268.48
$ws.Range("D7").Value = 10152.74
$ws.Range("E7").Value = 313.269989
$ws.Range("F7").Value = 309.700012
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = -1.14
$ws.Range("I7").Value = $true

# Match the date-style formatting used elsewhere in columns A and G (numFmtId 22)
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"
$ws.Range("G7").NumberFormat = "m/d/yy h:mm"

# Widen columns E and F to fit the new values
$ws.Columns("E:F").ColumnWidth = 10
